$wb = $excel.ActiveWorkbook

# Sheet 2: Главные
$ws = $wb.Worksheets.Item(2)
$ws.Range("AA2:AA26").Value = "2025-11-17 23:26:53"
# Row 4
$ws.Range("C4").Value = 19
$ws.Range("D4").Value = 312
$ws.Range("E4").Value = 131
$ws.Range("F4").Value = 181
$ws.Range("G4").Value = 16.42
$ws.Range("H4").Value = 6.89
$ws.Range("I4").Value = 9.529999999999999
$ws.Range("J4").Value = 63
$ws.Range("K4").Value = 78
# Row 14
$ws.Range("C14").Value = 18
$ws.Range("D14").Value = 220
$ws.Range("E14").Value = 115
$ws.Range("F14").Value = 105
$ws.Range("G14").Value = 12.22
$ws.Range("H14").Value = 6.39
$ws.Range("I14").Value = 5.83
$ws.Range("J14").Value = 55
$ws.Range("K14").Value = 45
$ws.Range("M14").Value = 3
# Row 15
$ws.Range("C15").Value = 17
$ws.Range("D15").Value = 290
$ws.Range("E15").Value = 137
$ws.Range("F15").Value = 153
$ws.Range("G15").Value = 17.06
$ws.Range("H15").Value = 8.06
$ws.Range("I15").Value = 9
$ws.Range("J15").Value = 51
$ws.Range("K15").Value = 69
$ws.Range("L15").Value = 3
$ws.Range("M15").Value = 3
$ws.Range("P15").Value = 1
$ws.Range("V15").Value = 4
# Row 22
$ws.Range("C22").Value = 19
$ws.Range("D22").Value = 410
$ws.Range("E22").Value = 181
$ws.Range("F22").Value = 229
$ws.Range("G22").Value = 21.58
$ws.Range("H22").Value = 9.529999999999999
$ws.Range("I22").Value = 12.05
$ws.Range("J22").Value = 73
$ws.Range("K22").Value = 77
$ws.Range("L22").Value = 3
$ws.Range("M22").Value = 3
$ws.Range("P22").Value = 1
$ws.Range("V22").Value = 6
# Row 25
$ws.Range("C25").Value = 26
$ws.Range("D25").Value = 429
$ws.Range("E25").Value = 208
$ws.Range("F25").Value = 221
$ws.Range("G25").Value = 16.5
$ws.Range("H25").Value = 8
$ws.Range("I25").Value = 8.5
$ws.Range("J25").Value = 99
$ws.Range("K25").Value = 103
$ws.Range("M25").Value = 3

# Sheet 3: Линейные
$ws = $wb.Worksheets.Item(3)
$ws.Range("AA2:AA26").Value = "2025-11-17 23:26:53"
# Row 8
$ws.Range("C8").Value = 23
$ws.Range("D8").Value = 365
$ws.Range("E8").Value = 136
$ws.Range("F8").Value = 229
$ws.Range("G8").Value = 15.87
$ws.Range("I8").Value = 9.960000000000001
$ws.Range("J8").Value = 63
$ws.Range("K8").Value = 87
# Row 17
$ws.Range("C17").Value = 13
$ws.Range("D17").Value = 240
$ws.Range("E17").Value = 141
$ws.Range("F17").Value = 99
$ws.Range("G17").Value = 18.46
$ws.Range("H17").Value = 10.85
$ws.Range("I17").Value = 7.62
$ws.Range("J17").Value = 53
$ws.Range("K17").Value = 42
$ws.Range("L17").Value = 3
$ws.Range("M17").Value = 3
$ws.Range("P17").Value = 1
$ws.Range("V17").Value = 8
# Row 21
$ws.Range("C21").Value = 28
$ws.Range("D21").Value = 553
$ws.Range("E21").Value = 232
$ws.Range("F21").Value = 321
$ws.Range("G21").Value = 19.75
$ws.Range("H21").Value = 8.289999999999999
$ws.Range("I21").Value = 11.46
$ws.Range("J21").Value = 106
$ws.Range("K21").Value = 128
$ws.Range("M21").Value = 3
# Row 22
$ws.Range("C22").Value = 19
$ws.Range("D22").Value = 363
$ws.Range("E22").Value = 187
$ws.Range("F22").Value = 176
$ws.Range("G22").Value = 19.11
$ws.Range("H22").Value = 9.84
$ws.Range("I22").Value = 9.26
$ws.Range("J22").Value = 76
$ws.Range("K22").Value = 78
$ws.Range("L22").Value = 3
$ws.Range("M22").Value = 4
$ws.Range("P22").Value = 1
$ws.Range("V22").Value = 16
# Row 24
$ws.Range("C24").Value = 27
$ws.Range("D24").Value = 484
$ws.Range("E24").Value = 193
$ws.Range("F24").Value = 291
$ws.Range("G24").Value = 17.93
$ws.Range("H24").Value = 7.15
$ws.Range("I24").Value = 10.78
$ws.Range("J24").Value = 84
$ws.Range("K24").Value = 108
$ws.Range("M24").Value = 5
# Row 26
$ws.Range("C26").Value = 24
$ws.Range("D26").Value = 477
$ws.Range("E26").Value = 205
$ws.Range("F26").Value = 272
$ws.Range("G26").Value = 19.88
$ws.Range("H26").Value = 8.539999999999999
$ws.Range("I26").Value = 11.33
$ws.Range("J26").Value = 80
$ws.Range("K26").Value = 81
